# This script reproduces the diff: the first data row (row 2, the
# "合肥·原&铁&崩 only展" event on 2024-03-23) was removed from both the
# "展览" sheet and the "全部类型" sheet, causing every subsequent event
# row to shift up by one position. In addition, several "想去人数" (F
# column) counters were refreshed to newer values as part of the same
# data sync.

$wb = $excel.ActiveWorkbook

# F-column ("想去人数") corrections to apply after the shift, keyed by
# the resulting row number once row 2 has been removed.
$fCorrections = @{
    2  = 7406
    3  = 7347
    9  = 93
    10 = 127
    12 = 94
    13 = 665
    14 = 505
    19 = 73
}

# Both the "展览" sheet (index 1) and the "全部类型" sheet (index 4)
# contain the same table and received the identical edit.
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Remove the old first event row; Excel shifts rows 3:20 up to 2:19
    # and the sheet dimension is recalculated automatically.
    $ws.Rows.Item(2).Delete()

    # The leftmost index column (A) is a plain 1..18 sequence that is
    # independent of which event occupies the row, so restore it after
    # the shift (Excel's row delete would otherwise have shifted those
    # numbers up along with everything else).
    for ($r = 2; $r -le 19; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply the refreshed "想去人数" counts.
    foreach ($row in $fCorrections.Keys) {
        $ws.Cells.Item($row, 6).Value = $fCorrections[$row]
    }
}
